# CHANGE Miilon - Add more vocabulary.
#
# 1) Add 4 new entries (no "grammar" value) to the end of the "LA PRESSE" sheet.
# 2) Add a brand-new sheet "LES TYPES D'ARTICLES" (after "LA PRESSE") with a
#    header row plus 9 new vocabulary entries, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: builds the same " { ""foreign"": ""...""... }," formula the sheet
# already uses for column F, anchored at a given row.
# ---------------------------------------------------------------------------
function Set-JsonFormula($ws, [int]$row) {
    $formula = ' "{ ""foreign"": """ & A' + $row + ' & """, ""grammar"": """ & B' + $row + ' & """, ""pronunciation"": """ & C' + $row + ' & """, ""meaning"": """ & D' + $row + ' & """ },"'
    $ws.Cells.Item($row, 6).Formula = "=" + $formula
}

# ---------------------------------------------------------------------------
# 1) "LA PRESSE" sheet - append rows 19-22
# ---------------------------------------------------------------------------
$press = $wb.Worksheets.Item("LA PRESSE")

$press.Cells.Item(19, 1).Value = "avoir bonne presse"
$press.Cells.Item(19, 3).Value = "avu^a:r bon pres"
$press.Cells.Item(19, 4).Value = "mít dobrou kritiku"
Set-JsonFormula $press 19

$press.Cells.Item(20, 1).Value = "une feuille de chou"
$press.Cells.Item(20, 3).Value = "ün föj d@ šu"
$press.Cells.Item(20, 4).Value = "zelný list; plátek (noviny)"
Set-JsonFormula $press 20

$press.Cells.Item(21, 1).Value = "la rubrique des chiens écrasés"
$press.Cells.Item(21, 3).Value = "rübrik de. šje~ e.kraze."
$press.Cells.Item(21, 4).Value = "novinová rubrika s příspěvky na různá nepříliš důležitá témata"
Set-JsonFormula $press 21

$press.Cells.Item(22, 1).Value = "un torchon"
$press.Cells.Item(22, 3).Value = "ö~ toršo~"
$press.Cells.Item(22, 4).Value = "hadr, bezcenný spis"
Set-JsonFormula $press 22

$press.Range("F19").Select()

# ---------------------------------------------------------------------------
# 2) New sheet "LES TYPES D'ARTICLES" placed right after "LA PRESSE"
# ---------------------------------------------------------------------------
$types = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $press)
$types.Name = "LES TYPES D'ARTICLES"

# Column widths to match the other sheets
$types.Columns.Item(1).ColumnWidth = $press.Columns.Item(1).ColumnWidth
$types.Columns.Item(2).ColumnWidth = $press.Columns.Item(2).ColumnWidth
$types.Columns.Item(3).ColumnWidth = $press.Columns.Item(3).ColumnWidth
$types.Columns.Item(4).ColumnWidth = $press.Columns.Item(4).ColumnWidth

# Header row, copied from the other sheets
$types.Cells.Item(1, 1).Value = "Mot français"
$types.Cells.Item(1, 2).Value = "Grammaire"
$types.Cells.Item(1, 3).Value = "Prononciation"
$types.Cells.Item(1, 4).Value = "Signification en tchèque"
$press.Range("A1:D1").Copy()
$types.Range("A1:D1").PasteSpecial(-4122)

# Data rows 2-10
$types.Cells.Item(2, 1).Value = "billet"
$types.Cells.Item(2, 2).Value = "nm"
$types.Cells.Item(2, 3).Value = "bije"
$types.Cells.Item(2, 4).Value = "malý novinový článek (často polemický nebo satirický)"
Set-JsonFormula $types 2

$types.Cells.Item(3, 1).Value = "brève"
$types.Cells.Item(3, 2).Value = "nf"
$types.Cells.Item(3, 3).Value = "bre:v"
$types.Cells.Item(3, 4).Value = "zprávička (na poslední chvíli nebo nepříliš důležitá)"
Set-JsonFormula $types 3

$types.Cells.Item(4, 1).Value = "chronique"
$types.Cells.Item(4, 2).Value = "nf"
$types.Cells.Item(4, 3).Value = "kronik"
$types.Cells.Item(4, 4).Value = "kronika; rubrika (v novinách)"
Set-JsonFormula $types 4

$types.Cells.Item(5, 1).Value = "critique"
$types.Cells.Item(5, 2).Value = "nf"
$types.Cells.Item(5, 3).Value = "kritik"
$types.Cells.Item(5, 4).Value = "kritika, odborný posudek"
Set-JsonFormula $types 5

$types.Cells.Item(6, 1).Value = "dépêche"
$types.Cells.Item(6, 2).Value = "nf"
$types.Cells.Item(6, 3).Value = "de.peš"
$types.Cells.Item(6, 4).Value = "zpráva (stručná)"
Set-JsonFormula $types 6

$types.Cells.Item(7, 1).Value = "éditorial"
$press.Cells.Item(7, 1).Copy()
$types.Cells.Item(7, 1).PasteSpecial(-4122)
$types.Cells.Item(7, 2).Value = "nm"
$types.Cells.Item(7, 3).Value = "e.ditorjal"
$types.Cells.Item(7, 4).Value = "editorial, úvodník"
Set-JsonFormula $types 7

$types.Cells.Item(8, 1).Value = "enquête"
$types.Cells.Item(8, 2).Value = "nf"
$types.Cells.Item(8, 3).Value = "a~ket"
$types.Cells.Item(8, 4).Value = "průzkum, anketa"
Set-JsonFormula $types 8

$types.Cells.Item(9, 1).Value = "entretien"
$types.Cells.Item(9, 2).Value = "nm"
$types.Cells.Item(9, 3).Value = "a~tr@tje~"
$types.Cells.Item(9, 4).Value = "rozhovor"
Set-JsonFormula $types 9

$types.Cells.Item(10, 1).Value = "reportage"
$types.Cells.Item(10, 2).Value = "nm"
$types.Cells.Item(10, 3).Value = "r@porta:ž"
$types.Cells.Item(10, 4).Value = "reportáž"
Set-JsonFormula $types 10

$types.Range("A1").Select()

# Make the new sheet the active one (matches activeTab="2" in workbook.xml)
$types.Activate()
